$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9042705297470093
$ws.Range("B1").Value = 1.094759225845337
$ws.Range("C1").Value = 0.8487009406089783
$ws.Range("D1").Value = 3.17652416229248
$ws.Range("E1").Value = 2.93181300163269
